$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    # Force the cell to accept the value as literal text even when it
    # looks like a number (e.g. "1.008"), without leaving a lasting
    # number-format / style change behind on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 4) '30.553.26'
Set-TextValue $ws.Cells.Item(2, 5) '  +0.44%  '
Set-TextValue $ws.Cells.Item(3, 4) '2.137.86'
Set-TextValue $ws.Cells.Item(3, 5) '  +1.83%  '
Set-TextValue $ws.Cells.Item(4, 4) '1.008'
Set-TextValue $ws.Cells.Item(4, 5) '  +0.38%  '
Set-TextValue $ws.Cells.Item(5, 4) '352.24'
Set-TextValue $ws.Cells.Item(5, 5) '  +5.47%  '
Set-TextValue $ws.Cells.Item(6, 4) '1.007'
Set-TextValue $ws.Cells.Item(6, 5) '  +0.45%  '
Set-TextValue $ws.Cells.Item(7, 4) '0.5261'
Set-TextValue $ws.Cells.Item(7, 5) '  +1.09%  '
Set-TextValue $ws.Cells.Item(8, 4) '0.4557'
Set-TextValue $ws.Cells.Item(8, 5) '  +0.49%  '
Set-TextValue $ws.Cells.Item(9, 4) '53.56'
Set-TextValue $ws.Cells.Item(9, 5) '  -1.52%  '
Set-TextValue $ws.Cells.Item(10, 4) '0.09179'
Set-TextValue $ws.Cells.Item(10, 5) '  +3.40%  '
Set-TextValue $ws.Cells.Item(11, 4) '1.185'
Set-TextValue $ws.Cells.Item(11, 5) '  +0.57%  '
Set-TextValue $ws.Cells.Item(12, 4) '25.43'
Set-TextValue $ws.Cells.Item(12, 5) '  +5.70%  '
Set-TextValue $ws.Cells.Item(13, 4) '2.139.46'
Set-TextValue $ws.Cells.Item(13, 5) '  +2.22%  '
Set-TextValue $ws.Cells.Item(14, 4) '6.906'
Set-TextValue $ws.Cells.Item(14, 5) '  +1.63%  '
Set-TextValue $ws.Cells.Item(15, 4) '8.205'
Set-TextValue $ws.Cells.Item(16, 4) '102.31'
Set-TextValue $ws.Cells.Item(16, 5) '  +5.65%  '
Set-TextValue $ws.Cells.Item(17, 4) '0.00001174'
Set-TextValue $ws.Cells.Item(17, 5) '  +2.85%  '
Set-TextValue $ws.Cells.Item(18, 4) '1.008'
Set-TextValue $ws.Cells.Item(18, 5) '  +0.49%  '
Set-TextValue $ws.Cells.Item(19, 4) '0.06720'
Set-TextValue $ws.Cells.Item(19, 5) '  +1.56%  '
Set-TextValue $ws.Cells.Item(20, 4) '20.44'
Set-TextValue $ws.Cells.Item(20, 5) '  +6.67%  '
Set-TextValue $ws.Cells.Item(21, 5) '  +0.44%  '
Set-TextValue $ws.Cells.Item(22, 4) '6.383'
Set-TextValue $ws.Cells.Item(22, 5) '  +1.83%  '
Set-TextValue $ws.Cells.Item(23, 4) '30.636.07'
Set-TextValue $ws.Cells.Item(23, 5) '  +0.59%  '
Set-TextValue $ws.Cells.Item(24, 4) '12.92'
Set-TextValue $ws.Cells.Item(24, 5) '  +4.91%  '
Set-TextValue $ws.Cells.Item(25, 4) '2.392'
Set-TextValue $ws.Cells.Item(25, 5) '  +2.41%  '
Set-TextValue $ws.Cells.Item(26, 4) '2.388.94'
Set-TextValue $ws.Cells.Item(26, 5) '  +2.16%  '
Set-TextValue $ws.Cells.Item(27, 2) 'EthereumClassic'
Set-TextValue $ws.Cells.Item(27, 3) 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Cells.Item(27, 4) '22.52'
Set-TextValue $ws.Cells.Item(27, 5) '  +1.61%  '
Set-TextValue $ws.Cells.Item(28, 2) 'LidoDAOToken'
Set-TextValue $ws.Cells.Item(28, 3) 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Cells.Item(28, 4) '2.655'
Set-TextValue $ws.Cells.Item(28, 5) '  +5.67%  '
Set-TextValue $ws.Cells.Item(29, 4) '164.92'
Set-TextValue $ws.Cells.Item(29, 5) '  +1.23%  '
Set-TextValue $ws.Cells.Item(30, 4) '136.10'
Set-TextValue $ws.Cells.Item(30, 5) '  +2.40%  '
Set-TextValue $ws.Cells.Item(31, 4) '1.226'
Set-TextValue $ws.Cells.Item(31, 5) '  +2.01%  '
Set-TextValue $ws.Cells.Item(32, 5) '  +1.65%  '
Set-TextValue $ws.Cells.Item(33, 4) '1.688'
Set-TextValue $ws.Cells.Item(33, 5) '  +2.11%  '
Set-TextValue $ws.Cells.Item(34, 4) '6.410'
Set-TextValue $ws.Cells.Item(34, 5) '  +0.78%  '
Set-TextValue $ws.Cells.Item(35, 4) '4.037'
Set-TextValue $ws.Cells.Item(35, 5) '  +2.43%  '
Set-TextValue $ws.Cells.Item(36, 4) '6.162'
Set-TextValue $ws.Cells.Item(36, 5) '  +6.36%  '
Set-TextValue $ws.Cells.Item(37, 4) '10.52'
Set-TextValue $ws.Cells.Item(37, 5) '  +1.26%  '
Set-TextValue $ws.Cells.Item(38, 4) '0.02653'
Set-TextValue $ws.Cells.Item(38, 5) '  +3.20%  '
Set-TextValue $ws.Cells.Item(39, 4) '0.06976'
Set-TextValue $ws.Cells.Item(39, 5) '  +2.18%  '
Set-TextValue $ws.Cells.Item(40, 4) '0.2338'
Set-TextValue $ws.Cells.Item(40, 5) '  +1.77%  '
Set-TextValue $ws.Cells.Item(41, 4) '12.78'
Set-TextValue $ws.Cells.Item(41, 5) '  +0.48%  '
Set-TextValue $ws.Cells.Item(42, 4) '0.7002'
Set-TextValue $ws.Cells.Item(42, 5) '  +2.07%  '
Set-TextValue $ws.Cells.Item(43, 4) '1.279'
Set-TextValue $ws.Cells.Item(43, 5) '  +2.70%  '
Set-TextValue $ws.Cells.Item(44, 4) '14.91'
Set-TextValue $ws.Cells.Item(44, 5) '  +7.09%  '
Set-TextValue $ws.Cells.Item(45, 4) '2.360'
Set-TextValue $ws.Cells.Item(45, 5) '  +1.83%  '
Set-TextValue $ws.Cells.Item(46, 4) '0.6513'
Set-TextValue $ws.Cells.Item(46, 5) '  +2.73%  '
Set-TextValue $ws.Cells.Item(47, 4) '0.00000000370'
Set-TextValue $ws.Cells.Item(47, 5) '  +7.75%  '
Set-TextValue $ws.Cells.Item(48, 4) '3.754'
Set-TextValue $ws.Cells.Item(48, 5) '  +2.69%  '
Set-TextValue $ws.Cells.Item(49, 4) '1.251'
Set-TextValue $ws.Cells.Item(49, 5) '  +0.45%  '
Set-TextValue $ws.Cells.Item(50, 4) '84.05'
Set-TextValue $ws.Cells.Item(50, 5) '  +1.36%  '
Set-TextValue $ws.Cells.Item(51, 4) '0.07296'
Set-TextValue $ws.Cells.Item(51, 5) '  +2.53%  '
